$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1,  "77H052.79",  "Bottom", "2025-11-27 23:16:56"),
    @(12, "GLZ50RLL",   "Bottom", "2025-11-27 23:16:52"),
    @(1,  "77H052.79",  "Bottom", "2025-11-27 23:42:14"),
    @(12, "177806279",  "Bottom", "2025-11-27 23:42:11"),
    @(1,  "77H052.79",  "Bottom", "2025-11-27 23:47:57"),
    @(12, "77H052.79",  "Bottom", "2025-11-27 23:47:54"),
    @(1,  "77H-052.79", "Bottom", "2025-11-27 23:55:10")
)

$startRow = 19
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $bCell = $ws.Cells.Item($r, 2)
    $plate = $row[1]

    $ws.Cells.Item($r, 1).Value = $row[0]

    if ($plate -match '^[0-9]+$') {
        # Purely-numeric plate text (e.g. "177806279") would otherwise be
        # auto-coerced to a number by Excel; force text so it round-trips
        # as a string, then restore the default cell style.
        $bCell.NumberFormat = "@"
        $bCell.Value = $plate
        $bCell.Style = "Normal"
    } else {
        $bCell.Value = $plate
    }

    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
